$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: UCRN, INSIGHT REF, MATTER TYPE
$ws.Range("C2").Value = "S459/53671"
$ws.Range("D2").Value = "MRT1835"
$ws.Range("E2").Value = "Purchase TL"

# Update row 3: UCRN, INSIGHT REF (MATTER TYPE stays "New Build Purchase")
$ws.Range("C3").Value = "S459/53977"
$ws.Range("D3").Value = "NBT1893"

# Remove row 4 entirely (shrinks the used range to A1:E3)
$ws.Rows.Item(4).Delete()

# Update the selection to E11 as reflected in the sheet view
$ws.Range("E11").Select()
